# Slide 3 ("EPRIME procedura" intro slide): the instruction line that used
# to read plainly now highlights the verb "ignorować" in red and the noun
# at the end of the sentence changes from "bodźców" to "strzałek".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Locate the "Content Placeholder 2" shape that holds the instruction text
# ("Proszę zupełnie ignorować kierunek poniższych bodźców:").
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -like "Prosz*zupe*ignorowa*") {
            $targetShape = $sh
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Update the wording first (bodźców -> strzałek), then recolor the
# "ignorować" word red; re-coloring a sub-range naturally splits the run
# into three runs, matching the target formatting.
$tr.Text = "Proszę zupełnie ignorować kierunek poniższych strzałek:"

$highlight = $tr.Characters(17, 9)
$highlight.Font.Color.RGB = 255
